# Add 6 new emergency event rows (EV_2025_57 .. EV_2025_62) to the "events" sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("events")
$wsOriginalActive = $wb.ActiveSheet

$rows = @(
    @{ A="EV_2025_57"; B=45992; C=45991; D="Huánuco";      E="Ambo";         F="Ambo";            G="Vientos fuertes";  H="En proceso"; I="Viviendas afectadas"; J="Evaluación en curso"; K=-10.150138999999999; L=-76.141176999999999 },
    @{ A="EV_2025_58"; B=45992; C=45992; D="Huánuco";      E="Puerto Inca";  F="Codo del Pozuzo"; G="Sismo";            H="Activa";     I="En monitoreo";        J="Evaluación en curso"; K=-9.6343639999999997;  L=-75.466538999999997 },
    @{ A="EV_2025_59"; B=45993; C=45988; D="Huancavelica"; E="Tayacaja";     F="Huaribamba";      G="Vientos fuertes";  H="En proceso"; I="Viviendas afectadas"; J="Evaluación en curso"; K=-12.262100999999999; L=-74.908894000000004 },
    @{ A="EV_2025_60"; B=45993; C=45993; D="San Martín";   E="Lamas";        F="Lamas";           G="Lluvias intensas"; H="Activa";     I="En monitoreo";        J="Evaluación en curso"; K=-6.3990600000000004;  L=-76.526042000000004 },
    @{ A="EV_2025_61"; B=45993; C=45993; D="San Martín";   E="Bellavista";   F="Bellavista";      G="Lluvias intensas"; H="Activa";     I="En monitoreo";        J="Evaluación en curso"; K=-6.9910240000000003;  L=-76.600607999999994 },
    @{ A="EV_2025_62"; B=45993; C=45989; D="Huancavelica"; E="Tayacaja";     F="Pichos";          G="Vientos fuertes";  H="En proceso"; I="En monitoreo";        J="Evaluación en curso"; K=-12.181583;           L=-74.933282000000005 }
)

$startRow = 58
$dateFmtSrc = $ws.Cells.Item(57, 2)   # existing date-formatted cell to copy number format from

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $row = $rows[$i]

    $ws.Cells.Item($r, 1).Value = $row.A
    $ws.Cells.Item($r, 2).Value = $row.B
    $ws.Cells.Item($r, 3).Value = $row.C
    $ws.Cells.Item($r, 4).Value = $row.D
    $ws.Cells.Item($r, 5).Value = $row.E
    $ws.Cells.Item($r, 6).Value = $row.F
    $ws.Cells.Item($r, 7).Value = $row.G
    $ws.Cells.Item($r, 8).Value = $row.H
    $ws.Cells.Item($r, 9).Value = $row.I
    $ws.Cells.Item($r, 10).Value = $row.J
    $ws.Cells.Item($r, 11).Value = $row.K
    $ws.Cells.Item($r, 12).Value = $row.L
}

# Copy the existing date number format (not a literal format string) onto the
# new B/C cells so no new custom numFmt gets created -- re-uses the same style.
$dateFmtSrc.Copy()
$ws.Range($ws.Cells.Item($startRow, 2), $ws.Cells.Item($startRow + $rows.Count - 1, 3)).PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Update the saved selection for the "events" sheet.
$ws.Range("M59").Select()

# Restore the sheet that was active before we touched "events" so the
# workbook's active tab / tabSelected state is unchanged on save.
$wsOriginalActive.Select()
